# Edit LOB1019.xlsx worksheet per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix column definitions: col A (1) should only be 30.71 wide; col B (2) handled by its own <col>.
$ws.Columns.Item(1).ColumnWidth = 30.7109375

# Insert a new row at 13 (pushes old rows 13-24 down to 14-25), matching the
# row that splits "Docentes responsaveis:" (row 12, label-only) from the rest.
$ws.Rows.Item(13).Insert()

# The inserted row's A13 picked up a bold style with no content from the row
# above; the target layout has no A13 cell at all (row 13 only has B/C).
$ws.Range("A13").Clear()

# Give B13/C13 the normal-wrap (col B) / red-wrap (col C) formatting used
# throughout the sheet, then fill them with the content that used to
# (incorrectly) sit in row 10.
$ws.Range("B2").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C13").PasteSpecial(-4122)
$ws.Range("B13").Value = '8711623 - Denize Kalempa'
$ws.Range("C13").Value = '8711623 - Denize Kalempa'

# Row 10 "Objetivos:" previously (incorrectly) held the professor's name;
# replace with the actual Portuguese objectives text.
$ws.Range("B10").Value = 'Introduzir os conhecimentos básicos sobre estática e dinâmica de fluidos, oscilações, ondas mecânicas e leis da termodinâmica.'
$ws.Range("C10").Value = 'Introduzir os conhecimentos básicos sobre estática e dinâmica de fluidos, oscilações, ondas mecânicas e leis da termodinâmica.'

# Row 14 "Programa resumido:" previously (incorrectly) held "Semestral";
# replace with the actual Portuguese short-syllabus text.
$ws.Range("B14").Value = 'Estática e dinâmica de fluidos, oscilações e ondas mecânicas, gases ideais, temperatura, calor e leis da termodinâmica.'
$ws.Range("C14").Value = 'Estática e dinâmica de fluidos, oscilações e ondas mecânicas, gases ideais, temperatura, calor e leis da termodinâmica.'

# Row 16 "Programa:" previously (incorrectly) held a date; replace with the
# actual Portuguese syllabus text.
$ws.Range("B16").Value = '1) Estática de fluidos: pressão, princípios de Pascal e Arquimedes, tensão superficial, capilaridade;2) Dinâmica de fluidos: vazão, fluidos ideais, equação da continuidade, equação de Bernoulli, viscosidade, lei de Hagen-Poiseuille;3) Oscilações: movimento harmônico simples, amortecido e forçado, ressonância; 4) Ondas: transversais e longitudinais, equação de onda, superposição, interferência, ondas estacionárias e ressonância, ondas sonoras, intensidade e nível sonoro, batimentos, efeito Doppler;5) Temperatura e calor: conceitos, escalas de temperatura, a lei zero da termodinâmica, dilatação térmica, absorção de calor por sólidos e líquidos, calor e trabalho, mecanismos de transferência de calor, gases ideais, calor específico molar de um gás ideal e graus de liberdade;6) Termodinâmica: primeira lei da termodinâmica, processos reversíveis eirreversíveis, entropia, segunda lei da termodinâmica, máquinas térmicas eeficiência.'
$ws.Range("C16").Value = '1) Estática de fluidos: pressão, princípios de Pascal e Arquimedes, tensão superficial, capilaridade;2) Dinâmica de fluidos: vazão, fluidos ideais, equação da continuidade, equação de Bernoulli, viscosidade, lei de Hagen-Poiseuille;3) Oscilações: movimento harmônico simples, amortecido e forçado, ressonância; 4) Ondas: transversais e longitudinais, equação de onda, superposição, interferência, ondas estacionárias e ressonância, ondas sonoras, intensidade e nível sonoro, batimentos, efeito Doppler;5) Temperatura e calor: conceitos, escalas de temperatura, a lei zero da termodinâmica, dilatação térmica, absorção de calor por sólidos e líquidos, calor e trabalho, mecanismos de transferência de calor, gases ideais, calor específico molar de um gás ideal e graus de liberdade;6) Termodinâmica: primeira lei da termodinâmica, processos reversíveis eirreversíveis, entropia, segunda lei da termodinâmica, máquinas térmicas eeficiência.'

# Row 19 "Metodo:" previously (incorrectly) held the professor's name;
# replace with the grading-method text.
$ws.Range("B19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'
$ws.Range("C19").Value = 'NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.'

# Row 20 "Criterio:" previously (incorrectly) held the grading-method text;
# replace with the passing-grade criterion.
$ws.Range("B20").Value = 'NF≥ 5,0.'
$ws.Range("C20").Value = 'NF≥ 5,0.'

# Row 21 "Norma de recuperacao:" previously (incorrectly) held the passing
# grade criterion; replace with the actual recovery norm text.
$ws.Range("B21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'
$ws.Range("C21").Value = '(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.'

# Row 22 "Bibliografia:" previously (incorrectly) held the recovery norm
# text; replace with the actual bibliography.
$ws.Range("B22").Value = 'NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 2, Edgard Blucher (2008).RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol.2, LTC (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.2, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 2, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 2, Thomson Pioneira (2008).'
$ws.Range("C22").Value = 'NUSSENZVEIG, H.M. Curso de Física Básica. Vol. 2, Edgard Blucher (2008).RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol.2, LTC (2008).TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol.2, LTC (2008).SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 2, Pearson Addison Wesley (2009).JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 2, Thomson Pioneira (2008).'
